$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New row 4 (Scenario 3 / delete endpoint): start from a copy of row 3 so
#    that formatting (styles) and any re-used text (content-Type, description,
#    name/phone/address/types/website/language, origin) line up exactly with
#    the existing rows, then overwrite just the cells that actually differ.
# ---------------------------------------------------------------------------
$ws.Range("A3:P3").Copy($ws.Range("A4:P4"))
$ws.Range("F4:G4").ClearContents()

# ---------------------------------------------------------------------------
# 2) New column Q (place_id): seed formatting from column P on each existing
#    row, then clear the row 2 cell back out (it stays blank in the diff).
# ---------------------------------------------------------------------------
$ws.Range("P1").Copy($ws.Range("Q1"))
$ws.Range("P2").Copy($ws.Range("Q2"))
$ws.Range("P3").Copy($ws.Range("Q3"))
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q2").ClearContents()

# ---------------------------------------------------------------------------
# 3) Fill in the actual new values. Order matters here: it reproduces the
#    shared-string table order seen in the target workbook.
# ---------------------------------------------------------------------------
$ws.Range("Q1").Value2 = "place_id"
$ws.Range("Q3").Value2 = "d6790ea8c04cea36517edcec20da0212"
$ws.Range("Q4").Value2 = "ef6ed47dffcf1a24b70ef776662f2bc1"
$ws.Range("A4").Value2 = "Scenario 3"
$ws.Range("D4").Value2 = "/maps/api/place/delete/json"

# ---------------------------------------------------------------------------
# 4) New hyperlink for E4, pointing at the same host used by E2/E3.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E4"), "https://rahulshettyacademy.com/")

# ---------------------------------------------------------------------------
# 5) Column Q width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(17).ColumnWidth = 33.85

# ---------------------------------------------------------------------------
# 6) View state: selection moves to Q2 (new column).
# ---------------------------------------------------------------------------
$ws.Range("Q2").Select()
